# Rotate the four reference URLs on the "References" slides: each
# paragraph's URL is replaced by the URL that used to sit in the next
# paragraph, with the last one wrapping back to the first.
$p = $ppt.ActivePresentation

$urlMap = @{
    "https://en.wikipedia.org/wiki/Main_Page" = "https://www.nih.gov/"
    "https://www.nih.gov/"                    = "https://scholar.google.com/"
    "https://scholar.google.com/"             = "https://www.jstor.org/"
    "https://www.jstor.org/"                  = "https://en.wikipedia.org/wiki/Main_Page"
}

$slideIndexes = @(7, 13, 19)

foreach ($slideIdx in $slideIndexes) {
    $s = $p.Slides.Item($slideIdx)
    foreach ($shape in $s.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count
            for ($i = 1; $i -le $paraCount; $i++) {
                $para = $tr.Paragraphs($i, 1)
                # Paragraphs().Text includes a trailing paragraph-mark (CR);
                # strip it before comparing against the plain URL strings.
                $oldText = $para.Text.TrimEnd("`r", "`n")
                if ($urlMap.ContainsKey($oldText)) {
                    $run = $para.Runs(1, 1)
                    $run.Text = $urlMap[$oldText]
                }
            }
        }
    }
}
